# Final Project Parabolic Calculations.xlsx - fix/rework
#
# Commit message: "Xfinal calculations are shit and don't work"
#
# Changes applied:
#   1. Time (F2) changed from 0.5 to 2 -> ripples through the H/I/K/L
#      (X1/Y1/X2/Y2) formula columns automatically on recalculation.
#   2. The "Xf" columns (J and M, rows 3-10) had their broken formulas
#      cleared out - row 2's Xf formulas (J2/M2) are left alone.
#   3. Selection ends up on M2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the Time input value.
$ws.Range("F2").Value = 2

# 2. Clear the non-working "Xf" formulas for rows 3 through 10
#    (row 2's Xf formulas and the already-empty row 11/12 cells are untouched).
$ws.Range("J3:J10").ClearContents()
$ws.Range("M3:M10").ClearContents()

# 3. Leave the selection where the user ended up.
$ws.Range("M2").Select()
